# Apply "added hospital pay reservation" edit to the report workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Client bonus percentage -> reset to 0
$ws.Range("E2").Value = 0

# Company name placeholder text
$ws.Range("E6").Value = "string"

# Line item 1: name, quantity, price, agreed factory price, sum
$ws.Range("D9").Value = "string"
$ws.Range("E9").Value = 10
$ws.Range("F9").Value = 2320
$ws.Range("G9").Value = 2340
$ws.Range("H9").Value = 23200

# Totals
$ws.Range("H31").Value = 23200
$ws.Range("G32").Value = "0.0% скидка билан"
$ws.Range("H32").Value = 23200
$ws.Range("H33").Value = 25984
